# running analysis with new gRNAs
#
# A new gRNA target region "Tat " (distinct from the existing "Tat/Rev" row)
# is inserted into the EF50/EF90 summary table. This pushes the Tat/Rev,
# Vif and Vpr rows down by one and adds a brand-new Vpr row at the bottom,
# plus a handful of updated counts (Gag/Pol, LTR, Tat/Rev, Vif).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the new bottom row (10) -------------------------
# Duplicate the formatting of the last existing data row (9, "Vpr") onto the
# brand-new row 10 via copy/paste-format, so the new cells pick up the same
# bold/centered/bordered look without introducing new style entries.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# --- Step 2: shift the bottom three rows' data down by one -----------------
# (old row 9 "Vpr" -> new row 10, old row 8 "Vif" -> new row 9,
#  old row 7 "Tat/Rev" -> new row 8), writing the updated values directly.

# Row 10 <- was row 9 (Vpr)
$ws.Range("A10").Value = "Vpr"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0

# Row 9 <- was row 8 (Vif), with updated EF50 count
$ws.Range("A9").Value = "Vif"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0

# Row 8 <- was row 7 (Tat/Rev), with updated EF90 count
$ws.Range("A8").Value = "Tat/Rev"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0

# --- Step 3: write the new "Tat " row into the now-vacated row 7 -----------
$ws.Range("A7").Value = "Tat "
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

# --- Step 4: updated counts elsewhere from the re-run analysis -------------
$ws.Range("B3").Value = 57
$ws.Range("C3").Value = 32
$ws.Range("B4").Value = 148
$ws.Range("C4").Value = 89
